# presurvey.xlsx update
# The respondent previously listed as "okankit@gmail.com" in C3 (ankit sahu's
# row) has their e-mail corrected to "okankit1312@gmail.com", with a mailto
# hyperlink added on that cell (matching the style already used for the other
# e-mail addresses in column C), and the active selection left on E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the e-mail address for row 3 (ankit sahu)
$ws.Range("C3").Value = "okankit1312@gmail.com"

# Add a mailto hyperlink for the corrected address, then restore the
# standard "Hyperlink" cell style (Add() re-applies it, but can leave the
# style index slightly different from the rest of the column - normalize it
# back to the same named style used elsewhere).
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:okankit1312@gmail.com")
$ws.Range("C3").Style = "Hyperlink"

# Leave the selection where the author left it when they saved.
$ws.Range("E10").Select()
